# Generate Report for Handback
# Adds a new row (row 4) for file 41adc32b-f722-4a45-854a-c90efbab7637.md
# to the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # decimal BGR for RGB(100,149,237) == FF6495ED

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Duplicate row 3's layout into row 4, then fix up the values/styles.
$ws.Range("A3:G3").Copy($ws.Range("A4"))

$ws.Range("A4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.md"
$ws.Range("B4").Value = "e2e\41adc32b-f722-4a45-854a-c90efbab7637.md"
$ws.Range("C4").Value = ".md"
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "Handed back: in sync with en-US"
$ws.Range("G4").Value = "2016-09-02 00:54:48"

$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41adc32b7220a4a45854ac90efbab76370000000/e2e/41adc32b-f722-4a45-854a-c90efbab7637.md", "", "", "e2e\41adc32b-f722-4a45-854a-c90efbab7637.md")
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Font.Underline = $true
$ws.Range("B4").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3:P3").Copy($ws.Range("A4"))

$ws.Range("A4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "e2e"
$ws.Range("E4").Value = "ht"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.31b9df73daa9a7f0aa5a1d08daf9d551808cd66b.zh-cn.xlf"
$ws.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H4").Value = "2016-09-02 00:54:43"
$ws.Range("I4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.md"
$ws.Range("J4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.31b9df73daa9a7f0aa5a1d08daf9d551808cd66b.zh-cn.xlf"
$ws.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K4").Value = "2016-09-02 00:55:02"
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = "True"
$ws.Range("N4").ClearContents()
$ws.Range("O4").Value = "False"
$ws.Range("P4").ClearContents()

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41adc32b7220a4a45854ac90efbab76370000000/e2e/41adc32b-f722-4a45-854a-c90efbab7637.md", "", "", "41adc32b-f722-4a45-854a-c90efbab7637.md")
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Font.Underline = $true
$ws.Range("A4").Font.Color = $hyperlinkColor

$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/31b9df73daa9a7f0aa5a1d08daf9d551808cd66b/e2e/41adc32b-f722-4a45-854a-c90efbab7637.md", "", "", "41adc32b-f722-4a45-854a-c90efbab7637.md")
$ws.Range("I4").ClearFormats()
$ws.Range("I4").Font.Underline = $true
$ws.Range("I4").Font.Color = $hyperlinkColor

$ws.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3:P3").Copy($ws.Range("A4"))

$ws.Range("A4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "e2e"
$ws.Range("E4").Value = "ht"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.31b9df73daa9a7f0aa5a1d08daf9d551808cd66b.de-de.xlf"
$ws.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H4").Value = "2016-09-02 00:54:48"
$ws.Range("I4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.md"
$ws.Range("J4").Value = "41adc32b-f722-4a45-854a-c90efbab7637.31b9df73daa9a7f0aa5a1d08daf9d551808cd66b.de-de.xlf"
$ws.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K4").Value = "2016-09-02 00:55:17"
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = "True"
$ws.Range("N4").ClearContents()
$ws.Range("O4").Value = "False"
$ws.Range("P4").ClearContents()

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41adc32b7220a4a45854ac90efbab76370000000/e2e/41adc32b-f722-4a45-854a-c90efbab7637.md", "", "", "41adc32b-f722-4a45-854a-c90efbab7637.md")
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Font.Underline = $true
$ws.Range("A4").Font.Color = $hyperlinkColor

$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/31b9df73daa9a7f0aa5a1d08daf9d551808cd66b/e2e/41adc32b-f722-4a45-854a-c90efbab7637.md", "", "", "41adc32b-f722-4a45-854a-c90efbab7637.md")
$ws.Range("I4").ClearFormats()
$ws.Range("I4").Font.Underline = $true
$ws.Range("I4").Font.Color = $hyperlinkColor

$ws.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
